$d = $word.ActiveDocument
$d.Content.Find.Execute("15.03.2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "16.03.2023", 2)
